# ---------------------------------------------------------------------------
# Scheduled market-data refresh for the Diabolos Profits workbook.
#
# Updates the currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ
# and the derived Leve price / profit columns (H, I, J, K, L, M, N) for a handful
# of Leve rows across the ALC, ARM, BSM, CRP, CUL and LTW sheets, reflecting newly
# pulled Universalis market-board prices.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ==== Sheet: ALC ====
$ws = $wb.Worksheets.Item("ALC")

# Row 20: H20=1987.5, I20=1987.5, K20=1987.5, M20=-1757.5
$ws.Range("H20").Value = 1987.5
$ws.Range("I20").Value = 1987.5
$ws.Range("K20").Value = 1987.5
$ws.Range("M20").Value = -1757.5

# Row 28: H28=44506.566, J28=2190.3333, L28=2190.3333, N28=-3160.3333
$ws.Range("H28").Value = 44506.566
$ws.Range("J28").Value = 2190.3333
$ws.Range("L28").Value = 2190.3333
$ws.Range("N28").Value = -3160.3333

# Row 33: H33=30166.76, I33=34262.273, K33=34262.273, M33=-34033.273
$ws.Range("H33").Value = 30166.76
$ws.Range("I33").Value = 34262.273
$ws.Range("K33").Value = 34262.273
$ws.Range("M33").Value = -34033.273

# Row 35: H35=1987.5, I35=1987.5, K35=1987.5, M35=-1608.5
$ws.Range("H35").Value = 1987.5
$ws.Range("I35").Value = 1987.5
$ws.Range("K35").Value = 1987.5
$ws.Range("M35").Value = -1608.5

# Row 62: H62=1551430, I62=3436668, K62=3436668, M62=-3436044
$ws.Range("H62").Value = 1551430
$ws.Range("I62").Value = 3436668
$ws.Range("K62").Value = 3436668
$ws.Range("M62").Value = -3436044

# Row 65: H65=1551430, I65=3436668, K65=17183340, M65=-17180220
$ws.Range("H65").Value = 1551430
$ws.Range("I65").Value = 3436668
$ws.Range("K65").Value = 17183340
$ws.Range("M65").Value = -17180220

# Row 111: H111=20139.5, I111=1451.375, J111=57515.75, K111=4354.125, L111=172547.25, M111=-1287.125, N111=-178681.25
$ws.Range("H111").Value = 20139.5
$ws.Range("I111").Value = 1451.375
$ws.Range("J111").Value = 57515.75
$ws.Range("K111").Value = 4354.125
$ws.Range("L111").Value = 172547.25
$ws.Range("M111").Value = -1287.125
$ws.Range("N111").Value = -178681.25

# Row 127: H127=1890.3572, I127=1331.5, J127=2635.5, K127=3994.5, L127=7906.5, M127=965.5, N127=-17826.5
$ws.Range("H127").Value = 1890.3572
$ws.Range("I127").Value = 1331.5
$ws.Range("J127").Value = 2635.5
$ws.Range("K127").Value = 3994.5
$ws.Range("L127").Value = 7906.5
$ws.Range("M127").Value = 965.5
$ws.Range("N127").Value = -17826.5

# Row 131: H131=35266.332, J131=35266.332, L131=105798.996, N131=-115878.996
$ws.Range("H131").Value = 35266.332
$ws.Range("J131").Value = 35266.332
$ws.Range("L131").Value = 105798.996
$ws.Range("N131").Value = -115878.996


# ==== Sheet: ARM ====
$ws = $wb.Worksheets.Item("ARM")

# Row 32: H32=1341.8889, I32=1341.8889, K32=1341.8889, M32=-1054.8889
$ws.Range("H32").Value = 1341.8889
$ws.Range("I32").Value = 1341.8889
$ws.Range("K32").Value = 1341.8889
$ws.Range("M32").Value = -1054.8889

# Row 61: H61=13890277, I61=13890277, K61=13890277, M61=-13890065
$ws.Range("H61").Value = 13890277
$ws.Range("I61").Value = 13890277
$ws.Range("K61").Value = 13890277
$ws.Range("M61").Value = -13890065

# Row 97: H97=862.5789, I97=494, K97=494, M97=2
$ws.Range("H97").Value = 862.5789
$ws.Range("I97").Value = 494
$ws.Range("K97").Value = 494
$ws.Range("M97").Value = 2

# Row 108: H108=64901, J108=64901, L108=64901, N108=-72581
$ws.Range("H108").Value = 64901
$ws.Range("J108").Value = 64901
$ws.Range("L108").Value = 64901
$ws.Range("N108").Value = -72581

# Row 122: H122=23813090, I122=30306388, J122=4326.3335, K122=90919164, L122=12979.0005, M122=-90916714, N122=-17879.0005
$ws.Range("H122").Value = 23813090
$ws.Range("I122").Value = 30306388
$ws.Range("J122").Value = 4326.3335
$ws.Range("K122").Value = 90919164
$ws.Range("L122").Value = 12979.0005
$ws.Range("M122").Value = -90916714
$ws.Range("N122").Value = -17879.0005

# Row 136: H136=13890277, I136=13890277, K136=41670831, M136=-41668281
$ws.Range("H136").Value = 13890277
$ws.Range("I136").Value = 13890277
$ws.Range("K136").Value = 41670831
$ws.Range("M136").Value = -41668281


# ==== Sheet: BSM ====
$ws = $wb.Worksheets.Item("BSM")

# Row 94: H94=2688.516, I94=2346.074, K94=2346.074, M94=-1895.074
$ws.Range("H94").Value = 2688.516
$ws.Range("I94").Value = 2346.074
$ws.Range("K94").Value = 2346.074
$ws.Range("M94").Value = -1895.074

# Row 99: H99=1270, I99=1200, K99=1200, M99=298
$ws.Range("H99").Value = 1270
$ws.Range("I99").Value = 1200
$ws.Range("K99").Value = 1200
$ws.Range("M99").Value = 298

# Row 102: H102=555.5, I102=555.5, K102=555.5, M102=2689.5
$ws.Range("H102").Value = 555.5
$ws.Range("I102").Value = 555.5
$ws.Range("K102").Value = 555.5
$ws.Range("M102").Value = 2689.5

# Row 107: H107=16669384, I107=2905.3845, K107=2905.3845, M107=-985.3845000000001
$ws.Range("H107").Value = 16669384
$ws.Range("I107").Value = 2905.3845
$ws.Range("K107").Value = 2905.3845
$ws.Range("M107").Value = -985.3845000000001


# ==== Sheet: CRP ====
$ws = $wb.Worksheets.Item("CRP")

# Row 16: H16=2274, I16=1365.3334, K16=1365.3334, M16=-1078.3334
$ws.Range("H16").Value = 2274
$ws.Range("I16").Value = 1365.3334
$ws.Range("K16").Value = 1365.3334
$ws.Range("M16").Value = -1078.3334

# Row 107: H107=2221.0833, I107=2361.4666, J107=1987.1111, K107=2361.4666, L107=1987.1111, M107=-441.4666000000002, N107=-5827.1111
$ws.Range("H107").Value = 2221.0833
$ws.Range("I107").Value = 2361.4666
$ws.Range("J107").Value = 1987.1111
$ws.Range("K107").Value = 2361.4666
$ws.Range("L107").Value = 1987.1111
$ws.Range("M107").Value = -441.4666000000002
$ws.Range("N107").Value = -5827.1111

# Row 113: H113=2274, I113=1365.3334, K113=1365.3334, M113=804.6666
$ws.Range("H113").Value = 2274
$ws.Range("I113").Value = 1365.3334
$ws.Range("K113").Value = 1365.3334
$ws.Range("M113").Value = 804.6666

# Row 122: H122=2956.6, I122=2876.375, K122=8629.125, M122=-6179.125
$ws.Range("H122").Value = 2956.6
$ws.Range("I122").Value = 2876.375
$ws.Range("K122").Value = 8629.125
$ws.Range("M122").Value = -6179.125

# Row 141: H141=110920, J141=110920, L141=110920, N141=-121280
$ws.Range("H141").Value = 110920
$ws.Range("J141").Value = 110920
$ws.Range("L141").Value = 110920
$ws.Range("N141").Value = -121280


# ==== Sheet: CUL ====
$ws = $wb.Worksheets.Item("CUL")

# Row 56: H56=6999.75, I56=6999.75, K56=6999.75, M56=-6469.75
$ws.Range("H56").Value = 6999.75
$ws.Range("I56").Value = 6999.75
$ws.Range("K56").Value = 6999.75
$ws.Range("M56").Value = -6469.75

# Row 68: H68=2777.7144, I68=2700, J68=2790.6667, K68=8100, L68=8372.000100000001, M68=-7289, N68=-9994.000100000001
$ws.Range("H68").Value = 2777.7144
$ws.Range("I68").Value = 2700
$ws.Range("J68").Value = 2790.6667
$ws.Range("K68").Value = 8100
$ws.Range("L68").Value = 8372.000100000001
$ws.Range("M68").Value = -7289
$ws.Range("N68").Value = -9994.000100000001

# Row 71: H71=2777.7144, I71=2700, J71=2790.6667, K71=24300, L71=25116.0003, M71=-20244, N71=-33228.0003
$ws.Range("H71").Value = 2777.7144
$ws.Range("I71").Value = 2700
$ws.Range("J71").Value = 2790.6667
$ws.Range("K71").Value = 24300
$ws.Range("L71").Value = 25116.0003
$ws.Range("M71").Value = -20244
$ws.Range("N71").Value = -33228.0003

# Row 98: H98=5208.8335, I98=14253.5, J98=686.5, K98=42760.5, L98=2059.5, M98=-41262.5, N98=-5055.5
$ws.Range("H98").Value = 5208.8335
$ws.Range("I98").Value = 14253.5
$ws.Range("J98").Value = 686.5
$ws.Range("K98").Value = 42760.5
$ws.Range("L98").Value = 2059.5
$ws.Range("M98").Value = -41262.5
$ws.Range("N98").Value = -5055.5


# ==== Sheet: LTW ====
$ws = $wb.Worksheets.Item("LTW")

# Row 40: H40=2583.1667, J40=2700, L40=2700, N40=-2972
$ws.Range("H40").Value = 2583.1667
$ws.Range("J40").Value = 2700
$ws.Range("L40").Value = 2700
$ws.Range("N40").Value = -2972

# Row 61: H61=14523.866, I61=10316.833, J61=17328.555, K61=10316.833, L61=17328.555, M61=-10114.833, N61=-17732.555
$ws.Range("H61").Value = 14523.866
$ws.Range("I61").Value = 10316.833
$ws.Range("J61").Value = 17328.555
$ws.Range("K61").Value = 10316.833
$ws.Range("L61").Value = 17328.555
$ws.Range("M61").Value = -10114.833
$ws.Range("N61").Value = -17732.555

# Row 82: H82=1476.2222, I82=1476.2222, K82=1476.2222, M82=-1115.2222
$ws.Range("H82").Value = 1476.2222
$ws.Range("I82").Value = 1476.2222
$ws.Range("K82").Value = 1476.2222
$ws.Range("M82").Value = -1115.2222

# Row 85: H85=1476.2222, I85=1476.2222, K85=1476.2222, M85=-228.2221999999999
$ws.Range("H85").Value = 1476.2222
$ws.Range("I85").Value = 1476.2222
$ws.Range("K85").Value = 1476.2222
$ws.Range("M85").Value = -228.2221999999999

# Row 100: H100=3124.4167, I100=2977, J100=3566.6667, K100=2977, L100=3566.6667, M100=-2436, N100=-4648.6667
$ws.Range("H100").Value = 3124.4167
$ws.Range("I100").Value = 2977
$ws.Range("J100").Value = 3566.6667
$ws.Range("K100").Value = 2977
$ws.Range("L100").Value = 3566.6667
$ws.Range("M100").Value = -2436
$ws.Range("N100").Value = -4648.6667

# Row 113: H113=14523.866, I113=10316.833, J113=17328.555, K113=10316.833, L113=17328.555, M113=-8146.833000000001, N113=-21668.555
$ws.Range("H113").Value = 14523.866
$ws.Range("I113").Value = 10316.833
$ws.Range("J113").Value = 17328.555
$ws.Range("K113").Value = 10316.833
$ws.Range("L113").Value = 17328.555
$ws.Range("M113").Value = -8146.833000000001
$ws.Range("N113").Value = -21668.555

# Row 117: H117=50000, J117=50000, L117=50000, N117=-59178
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -59178

# Row 122: H122=4273.077, I122=2935.7144, K122=8807.143199999999, M122=-6357.143199999999
$ws.Range("H122").Value = 4273.077
$ws.Range("I122").Value = 2935.7144
$ws.Range("K122").Value = 8807.143199999999
$ws.Range("M122").Value = -6357.143199999999

# Row 130: H130=0, J130=0, L130=0, N130=<cleared>
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 132: H132=3056.5144, I132=2701.3157, J132=3478.3125, K132=8103.9471, L132=10434.9375, M132=-5573.9471, N132=-15494.9375
$ws.Range("H132").Value = 3056.5144
$ws.Range("I132").Value = 2701.3157
$ws.Range("J132").Value = 3478.3125
$ws.Range("K132").Value = 8103.9471
$ws.Range("L132").Value = 10434.9375
$ws.Range("M132").Value = -5573.9471
$ws.Range("N132").Value = -15494.9375
